$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aciklama = "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."

$gorseller = @(
    "5010yelekkahverengi.jpg",
    "5010yelekkrem.jpg",
    "5010yeleksiyah.jpg",
    "5010yelektaş.jpg"
)

$urunAdlari = @(
    "ŞİŞME YELEK 5010 YELEK KAHVERENGİ",
    "ŞİŞME YELEK 5010 YELEK KREM",
    "ŞİŞME YELEK 5010 YELEK SİYAH",
    "ŞİŞME YELEK 5010 YELEKTAŞ"
)

$startRow = 85

# Column D (gorsel) filled first for all new rows
for ($i = 0; $i -lt $gorseller.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $gorseller[$i]
}

# Then column A (urun_adi) filled for all new rows
for ($i = 0; $i -lt $urunAdlari.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $urunAdlari[$i]
}

# Remaining columns (values already exist elsewhere as shared strings)
for ($i = 0; $i -lt 4; $i++) {
    $satir = $startRow + $i
    $ws.Cells.Item($satir, 2).Value = "500 TL"
    $ws.Cells.Item($satir, 3).Value = "Yelek"
    $ws.Cells.Item($satir, 5).Value = $aciklama
    $ws.Cells.Item($satir, 6).Value = "Var"
}

$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("H87").Select()
